$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B1").Formula = "=2.25*A1"
$ws.Range("B2:B48").Formula = "=2.25*A2"
$ws.Range("B49").Formula = "=2.25*A49"
$ws.Activate()
$ws.Range("F48").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
